$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.1345010810828807
$ws.Range("C2").Value = -0
$ws.Range("D2").Value = 0.2273454682876108
$ws.Range("E2").Value = 0.00242929335295269
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = -0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = -0.03466512907320284
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.2285361925391075
$ws.Range("N2").Value = 0.004469526511393119
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.1057514036729393
$ws.Range("V2").Value = 0.01736233589227211
$ws.Range("W2").Value = -0.04317109139228038
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = -0.07605427972340863
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = -0.01085648554341784
$ws.Range("AF2").Value = 0.005367920726140095
$ws.Range("AG2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = -0
$ws.Range("AL2").Value = -0.04746778189235903
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.0298602511210372
$ws.Range("AO2").Value = 0.07191143125891197
$ws.Range("AR2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1717595909833013
$ws.Range("AW2").Value = 0.08284363933337956
$ws.Range("AX2").Value = 0.01865155208229204
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.04079309510024247
$ws.Range("BF2").Value = 0.066304162275427
$ws.Range("BG2").Value = 0.03283340620120182
$ws.Range("BI2").Value = -0
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.02284731819459275
$ws.Range("BO2").Value = -0.05095449535107491
$ws.Range("BP2").Value = -0.06139456390466035
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.0224160439022575
$ws.Range("BX2").Value = 0.01880520505166924
$ws.Range("BY2").Value = -0.005576734370991781
$ws.Range("BZ2").Value = -0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = 0
$ws.Range("CE2").Value = 0.03627934320947346
$ws.Range("CG2").Value = -0.04186594810183301
$ws.Range("CH2").Value = 0.0295740607249829
$ws.Range("CJ2").Value = -0
$ws.Range("CL2").Value = -0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.02030754756250808
$ws.Range("CO2").Value = -0
$ws.Range("CP2").Value = 0.02985302072380105
$ws.Range("CQ2").Value = 0.05606994475932985
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04935248336591598
$ws.Range("CY2").Value = -0.04417562920233499
$ws.Range("CZ2").Value = 0.02296453756375039
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.04407648556797538
$ws.Range("DH2").Value = 0.01310168645478456
$ws.Range("DI2").Value = 0.05915278330528921
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DL2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.01774329046133111
$ws.Range("DQ2").Value = 0.03660075961852761
$ws.Range("DR2").Value = -0.01146582227784784
$ws.Range("DS2").Value = -0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.04865450827884053
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = -0.005347755318757407
$ws.Range("EA2").Value = -0.02853121857332327
$ws.Range("EB2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.04820678582620177
$ws.Range("EI2").Value = 0.08372349045933243
$ws.Range("EJ2").Value = -0.01623233054878248
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.04789466871976065
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.03644610419047922
$ws.Range("ES2").Value = 0.03649721248365179
$ws.Range("ET2").Value = 0
$ws.Range("EU2").Value = -0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.05197321496315279
$ws.Range("FA2").Value = -0.03525713401926399
$ws.Range("FB2").Value = 0.0279484242845618
$ws.Range("FD2").Value = -0
$ws.Range("FF2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = 0.01251368300357263
$ws.Range("FI2").Value = 0
$ws.Range("FJ2").Value = -0.02026803619954162
$ws.Range("FK2").Value = 0.001524862824296119
$ws.Range("FN2").Value = -0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = -0.0005064520923302226
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = -0.0009808121277767718
$ws.Range("FT2").Value = 0.0138916870285865
$ws.Range("FV2").Value = -0
$ws.Range("FW2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.04312959506008182
$ws.Range("GA2").Value = -0
$ws.Range("GB2").Value = 0.0299007786849073
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
